$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# about [C2]
$ws.Range('C2').Value = 'About %1$s! .🎄 <b>Android</b> .🤣 " @ ? ' + [char]9 + ' '' ©'

# biometric_lock [F28]
$ws.Range('F28').Value = 'Verrouiller l''application avec la biométrie ou le code PIN de l''appareil'

# biometrics_failure [F30]
$ws.Range('F30').Value = 'Échec de l''authentification biométrique/code PIN'

# biometrics_no_support [F31]
$ws.Range('F31').Value = 'Aucune fonction biométrique n''est disponible sur cet appareil'

# biometrics_not_setup [F32]
$ws.Range('F32').Value = 'Aucun verrouillage biométrique/code PIN n''est paramétré sur cet appareil'

# cant_add_files_PLURALS_one [F39]
$ws.Range('F39').Value = '%d fichier n''a pas pu être ajouté'

# cant_add_files_PLURALS_other [F40]
$ws.Range('F40').Value = '%d fichiers n''ont pas pu être ajoutés'

# cant_add_images_PLURALS_one [F45]
$ws.Range('F45').Value = 'Impossible d''ajouter %d image'

# cant_add_images_PLURALS_other [F46]
$ws.Range('F46').Value = 'Impossible d''ajouter %d images'

# cant_load_image [F52]
$ws.Range('F52').Value = 'Impossible de charger l''image. Elle a peut-être été déplacée ou supprimée.'

# choose_other_app [F60]
$ws.Range('F60').Value = 'Choisissez l''application depuis laquelle importer'

# content_density [F66]
$ws.Range('F66').Value = 'Densité d''affichage'

# create_new [F72]
$ws.Range('F72').Value = 'Il n''existe pas encore d''étiquettes. En créer une ?'

# delete_audio_recording_forever [F82]
$ws.Range('F82').Value = 'Supprimer définitivement l''enregistrement ?'

# delete_file [C84]
$ws.Range('C84').Value = 'Delete file ''%s''?'

# delete_file [D84]
$ws.Range('D84').Value = 'Datei ''%s'' löschen?'

# delete_file [F84]
$ws.Range('F84').Value = 'Supprimer le fichier ''%s''?'

# delete_image_forever [F86]
$ws.Range('F86').Value = 'Supprimer l''image définitivement ?'

# delete_label [F87]
$ws.Range('F87').Value = 'Supprimer l''étiquette ?'

# edit_label [F109]
$ws.Range('F109').Value = 'Modifier l''étiquette'

# error_while_renaming_image [F117]
$ws.Range('F117').Value = 'Erreur lors du renommage de l''image'

# evernote_help [C119]
$ws.Range('C119').Value = 'In order to import your Notes from Evernote you must export your Evernote Notebook as ENEX. Click Help to get more information.' + [char]10 + 'If you already have a ENEX file, click Import and choose it.'

# evernote_help [D119]
$ws.Range('D119').Value = 'Um deine Notizen von Evernote zu importierten, exportiere dein Evernote Notebook als ENEX. Klicke Hilfe für mehr Infos.' + [char]10 + 'Falls du schon ein ENEX Datei hast, klicke Import und wähle es aus.'

# evernote_help [F119]
$ws.Range('F119').Value = 'Pour importer vos notes depuis Evernote, vous devez exporter votre carnet de notes Evernote au format ENEX. Cliquez sur "Aide" pour plus d''information.' + [char]10 + 'Si vous possédez déjà un fichier ENEX, cliquez sur "Importer" et sélectionnez-le.'

# evernote_help [G119]
$ws.Range('G119').Value = 'Per importare le note da Evernote devi esportare il tuo Notebook Evernote come ENEX. Clicca su Aiuto per ulteriori informazioni.' + [char]10 + 'Se hai già un file ENEX, clicca Importa e selezionalo.'

# export_settings_failure [F123]
$ws.Range('F123').Value = 'Échec de l''exportation des paramètres, avez-vous choisi un chemin invalide ?'

# export_settings_message [C124]
$ws.Range('C124').Value = 'All Settings will be exported to a JSON file, which can be used to re-import stored settings.' + [char]10 + 'Be aware, that this does not include encrypted settings like the auto-backup password or the biometric encryption key'

# export_settings_message [D124]
$ws.Range('D124').Value = 'Alle Einstellungen werden als JSON Datei exportiert, die dann zum Re-Import genutzt werden kann.' + [char]10 + 'Beachte das hier bei verschlüsselte Einstellungen wie das Auto-Backup Passwort oder die biometrische Verschlüssung nicht exportiert wird'

# export_settings_message [F124]
$ws.Range('F124').Value = 'Tous les paramètres seront exportés dans un fichier JSON, qui pourra être utilisé pour réimporter les paramètres enregistrés.' + [char]10 + 'Veuillez noter que cela n''inclut pas les paramètres chiffrés, tels que le mot de passe de sauvegarde automatique ou la clé de chiffrement biométrique.'

# external_data_message [C127]
$ws.Range('C127').Value = 'By enabling this, the app’s internal database will be moved into the app’s external storage (Android/media/com.philkes.notallyx).' + [char]10 + 'In combination with file synchronization apps this can be used to synchronize NotallyX data between multiple devices.'

# external_data_message [F127]
$ws.Range('F127').Value = 'En activant cette option, la base de données interne de l''application sera déplacée vers le stockage externe de l''application (Android/media/com.philkes.notallyx).' + [char]10 + 'En combinaison avec des applications de synchronisation de fichiers, cela peut être utilisé pour synchroniser les données de NotallyX entre plusieurs appareils.'

# external_data_message [G127]
$ws.Range('G127').Value = 'Attivandolo, il database interno dell’app verrà spostato nella memoria esterna dell’app (Android/media/com.philkes.notallyx).' + [char]10 + 'Abbinandolo con un’app di sincronizzazione file potrai così sincronizzare i dati di NotallyX tra dispositivi diversi.'

# google_keep_help [C133]
$ws.Range('C133').Value = 'In order to import your Notes from Google Keep you must download your Google Takeout ZIP file. Click Help to get more information.' + [char]10 + 'If you already have a Takeout ZIP file, click Import and choose the ZIP file.'

# google_keep_help [D133]
$ws.Range('D133').Value = 'Um deine Notizen aus Google Notizen zu importieren musst du deine Google Takeout ZIP Datei herunterladen' + [char]10 + 'Falls du das Takeout ZIP schon hast, klicke auf Import und wähle es aus.'

# google_keep_help [F133]
$ws.Range('F133').Value = 'Pour importer vos notes depuis Google Keep, vous devez télécharger votre fichier ZIP Google Takeout. Cliquez sur "Aide" pour plus d' + [char]10 + [char]9 + 'information.' + [char]10 + 'Si vous possédez déjà un fichier ZIP Google Takeout, cliquez sur "Importer" et choisissez le fichier ZIP.'

# google_keep_help [G133]
$ws.Range('G133').Value = 'Per importare le note da Google Keep devi scaricare il tuo file ZIP di Google Takeout. Clicca su Aiuto per ulteriori informazioni.' + [char]10 + 'Se hai già un file ZIP di Takeout, clicca Importa e selezionalo.'

# image_format_not_supported [F136]
$ws.Range('F136').Value = 'Format d''image non supporté'

# import_backup_password_hint [F139]
$ws.Range('F139').Value = 'Si votre sauvegarde n''est pas protégée par mot de passe, cliquez seulement sur "Importer une sauvegarde", sinon entrez le mot de passe correspondant.'

# import_settings_failure [F142]
$ws.Range('F142').Value = 'Échec de l''importation des paramètres, avez-vous choisi le bon fichier ?'

# insert_an_sd_card_audio [F154]
$ws.Range('F154').Value = 'Insérer une carte SD pour enregistrer de l''audio'

# label_exists [F166]
$ws.Range('F166').Value = 'L''étiquette existe déjà'

# label_visibility [F167]
$ws.Range('F167').Value = 'Masquer/Afficher l''étiquette dans le panneau de navigation'

# labels_hidden_in_overview [F169]
$ws.Range('F169').Value = 'En activant cette option, les étiquettes des notes seront masquées dans la vue d''ensemble'

# labels_hidden_in_overview_title [F170]
$ws.Range('F170').Value = 'Masquer les étiquettes dans la vue d''ensemble'

# max_items_to_display [F181]
$ws.Range('F181').Value = 'Nombre maximum d''éléments à afficher pour une liste'

# max_labels_to_display [F182]
$ws.Range('F182').Value = 'Nombre maximum d''étiquettes à afficher dans le panneau de navigation'

# notes_will_be [C202]
$ws.Range('C202').Value = 'Notes will be backed up to your phone’s internal storage everyday.' + [char]10 + 'This may not work if you have power saving mode enabled'

# notes_will_be [D202]
$ws.Range('D202').Value = 'Notizen werden täglich auf dem internen Speicher deines Telefons gesichert.' + [char]10 + 'Dies funktioniert möglicherweise nicht, wenn du den Energiesparmodus aktiviert hast.'

# notes_will_be [E202]
$ws.Range('E202').Value = 'Se realizará una copia de seguridad de las notas en el almacenamiento interno de su teléfono todos los días.' + [char]10 + 'Es posible que esto no funcione si tiene habilitado el modo de ahorro de energía'

# notes_will_be [F202]
$ws.Range('F202').Value = 'Les notes seront sauvegardées dans le stockage interne de votre téléphone tous les jours.' + [char]10 + 'Cela peut ne pas fonctionner si le mode d''économie d''énergie est activé'

# notes_will_be [G202]
$ws.Range('G202').Value = 'Il backup delle note verrà salvato nella memoria interna del tuo telefono ogni giorno.' + [char]10 + 'Ciò potrebbe non funzionare se hai la modalità di risparmio energetico attivata'

# notes_will_be [I202]
$ws.Range('I202').Value = '笔记每天将备份到你手机的内部存储中。' + [char]10 + '如果你启用了节电模式，此功能可能无法正常工作。'

# plain_text_files_help [D211]
$ws.Range('D211').Value = 'Um deine Text-Notizen (einzele Datei oder Ordner) zu importieren, klicke Import.' + [char]10 + 'Jede Datei wird als einzelne Notiz importiert, der Dateiname wird zum Notiz-Titel. Sollte der Textinhalt mit einer List-Syntax beginnen (z.B. Markdown ’- [x]’, NotallyX syntax ’[✓]’, or ’*’, ’-’), wird die Datei als List-Notiz importiert.'

# please_grant_notally_audio [F213]
$ws.Range('F213').Value = 'Veuillez accorder à NotallyX l''autorisation d''enregistrer de l''audio. Les enregistrements ne quittent jamais votre appareil.'

# please_grant_notally_notification [F214]
$ws.Range('F214').Value = ' Veuillez accordez à NotallyX l''autorisation d''envoyer des notifications. Cela permet d''afficher la progression des opérations telles que la suppression d''images ou l''importation de sauvegarde si elles prennent du temps.'

# save_recording [F235]
$ws.Range('F235').Value = 'Sauvegarder l''enregistrement ?'

# save_to_device [F236]
$ws.Range('F236').Value = 'Sauvegarder sur l''appareil'

# saved_to_device [F237]
$ws.Range('F237').Value = 'Sauvegardé sur l''appareil'

# something_went_wrong_audio [C251]
$ws.Range('C251').Value = 'Something went wrong. The audio recording may have been moved or deleted.' + [char]10 + 'Error : (%1$d, %2$d)'

# something_went_wrong_audio [D251]
$ws.Range('D251').Value = 'Etwas lief schief. Die Sprachnotiz wurde eventuell verschoben oder gelöscht.' + [char]10 + 'Fehler : (%1$d, %2$d)'

# something_went_wrong_audio [F251]
$ws.Range('F251').Value = 'Une erreur est survenue. l''enregistrement audio a peut-être été déplacé ou supprimé.' + [char]10 + 'Erreur : (%1$d, %2$d)'

# something_went_wrong_audio [G251]
$ws.Range('G251').Value = 'Qualcosa è andato storto. La registrazione audio potrebbe essere stata spostata o rimossa.' + [char]10 + 'Errore : (%1$d, %2$d)'

# tap_for_more_options [F258]
$ws.Range('F258').Value = 'Toucher pour afficher plus d''options'

# to_record_audio [F263]
$ws.Range('F263').Value = 'Pour enregistrer de l''audio, autorisez NotallyX à accéder à votre microphone. Cliquez sur Paramètres > Autorisations et activez le microphone.'
